$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "*"
$ws.Range("C2").Value = "*"
$ws.Range("C3").Value = "#"
$ws.Range("C4").Value = "#"

$ws.Range("C5").Select()
